$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 14.6

$values = @(
    0.020878416066977934,
    -0.0099999995783441875,
    -0.0089999995842138247,
    0.28399150439170739,
    -0.005999999599361594,
    -0.0059999995888198043,
    -0.01999999951572029,
    -0.019999999512559263,
    -0.0059999995815269713,
    -0.0059999995792097138,
    0.031593149921050667,
    -0.0059999995774031589,
    -0.0059999995707045173,
    -0.011999999538000239,
    -0.0059999995674866469,
    -0.0059999995662223249,
    -0.0059999995645290127,
    -0.0089999995487390905,
    -0.0089999995831679946,
    -0.0089999995796929966,
    -0.0089999995792213738,
    -0.0089999995788776488,
    -0.02430465675365312,
    -0.041999999403531163,
    -0.04199999940049004,
    -0.0059999995881447887,
    -0.020685454693705196,
    -0.0059999995847048737,
    -0.011999999552465113,
    -0.019999999511093769,
    -0.014999999537737452,
    0.032678656393295924,
    -0.0059999995845165799
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
